# Auto-generated edit script: updates market-data columns (H-N)
# on each sheet per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 101.333336
$ws.Range("I5").Value = 86
$ws.Range("J5").Value = 132
$ws.Range("K5").Value = 86
$ws.Range("L5").Value = 132
$ws.Range("M5").Value = 29
$ws.Range("N5").Value = -362
$ws.Range("H39").Value = 97
$ws.Range("I39").Value = 97
$ws.Range("K39").Value = 291
$ws.Range("M39").Value = 5
$ws.Range("H64").Value = 4723.75
$ws.Range("I64").Value = 4396.6665
$ws.Range("K64").Value = 4396.6665
$ws.Range("M64").Value = -4148.6665
$ws.Range("H67").Value = 4723.75
$ws.Range("I67").Value = 4396.6665
$ws.Range("K67").Value = 4396.6665
$ws.Range("M67").Value = -3538.6665
$ws.Range("H70").Value = 13217
$ws.Range("I70").Value = 1099.5
$ws.Range("J70").Value = 19275.75
$ws.Range("K70").Value = 3298.5
$ws.Range("L70").Value = 57827.25
$ws.Range("M70").Value = -3028.5
$ws.Range("N70").Value = -58367.25
$ws.Range("H73").Value = 13217
$ws.Range("I73").Value = 1099.5
$ws.Range("J73").Value = 19275.75
$ws.Range("K73").Value = 3298.5
$ws.Range("L73").Value = 57827.25
$ws.Range("M73").Value = -2362.5
$ws.Range("N73").Value = -59699.25
$ws.Range("H76").Value = 7977.5557
$ws.Range("I76").Value = 7899
$ws.Range("K76").Value = 7899
$ws.Range("M76").Value = -7584
$ws.Range("H79").Value = 7977.5557
$ws.Range("I79").Value = 7899
$ws.Range("K79").Value = 7899
$ws.Range("M79").Value = -6807
$ws.Range("H86").Value = 5848.9
$ws.Range("I86").Value = 5994.5
$ws.Range("K86").Value = 5994.5
$ws.Range("M86").Value = -4871.5
$ws.Range("H89").Value = 5848.9
$ws.Range("I89").Value = 5994.5
$ws.Range("K89").Value = 29972.5
$ws.Range("M89").Value = -24356.5
$ws.Range("H106").Value = 32190.375
$ws.Range("I106").Value = 36003.285
$ws.Range("K106").Value = 36003.285
$ws.Range("M106").Value = -35372.285
$ws.Range("H135").Value = 730.08
$ws.Range("I135").Value = 503.9
$ws.Range("J135").Value = 1634.8
$ws.Range("K135").Value = 4535.099999999999
$ws.Range("L135").Value = 14713.2
$ws.Range("M135").Value = -2000.099999999999
$ws.Range("N135").Value = -19783.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 4100
$ws.Range("I14").Value = 2300
$ws.Range("K14").Value = 2300
$ws.Range("M14").Value = -2125
$ws.Range("H32").Value = 14721.71
$ws.Range("J32").Value = 27526.84
$ws.Range("L32").Value = 27526.84
$ws.Range("N32").Value = -28100.84
$ws.Range("H45").Value = 1100.2858
$ws.Range("I45").Value = 825.5
$ws.Range("J45").Value = 1466.6666
$ws.Range("K45").Value = 825.5
$ws.Range("L45").Value = 1466.6666
$ws.Range("M45").Value = -448.5
$ws.Range("N45").Value = -2220.6666
$ws.Range("H110").Value = 6084
$ws.Range("I110").Value = 6639.3
$ws.Range("J110").Value = 4233
$ws.Range("K110").Value = 6639.3
$ws.Range("L110").Value = 4233
$ws.Range("M110").Value = -4594.3
$ws.Range("N110").Value = -8323
$ws.Range("H132").Value = 1959.5555
$ws.Range("I132").Value = 1913.1818
$ws.Range("K132").Value = 5739.5454
$ws.Range("M132").Value = -3209.5454
$ws.Range("H134").Value = 98000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1308.8334
$ws.Range("I86").Value = 1340.7
$ws.Range("J86").Value = 1149.5
$ws.Range("K86").Value = 1340.7
$ws.Range("L86").Value = 1149.5
$ws.Range("M86").Value = -217.7
$ws.Range("N86").Value = -3395.5
$ws.Range("H89").Value = 1308.8334
$ws.Range("I89").Value = 1340.7
$ws.Range("J89").Value = 1149.5
$ws.Range("K89").Value = 6703.5
$ws.Range("L89").Value = 5747.5
$ws.Range("M89").Value = -1087.5
$ws.Range("N89").Value = -16979.5
$ws.Range("H105").Value = 4327.4517
$ws.Range("I105").Value = 3832.4707
$ws.Range("J105").Value = 4928.5
$ws.Range("K105").Value = 3832.4707
$ws.Range("L105").Value = 4928.5
$ws.Range("M105").Value = -2085.4707
$ws.Range("N105").Value = -8422.5
$ws.Range("H134").Value = 1416.0667
$ws.Range("I134").Value = 711
$ws.Range("K134").Value = 2133
$ws.Range("M134").Value = 402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.3125
$ws.Range("I7").Value = 69
$ws.Range("J7").Value = 43.625
$ws.Range("K7").Value = 69
$ws.Range("L7").Value = 43.625
$ws.Range("M7").Value = 44
$ws.Range("N7").Value = -269.625
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H62").Value = 34860.848
$ws.Range("I62").Value = 4391
$ws.Range("J62").Value = 70409
$ws.Range("K62").Value = 4391
$ws.Range("L62").Value = 70409
$ws.Range("M62").Value = -3767
$ws.Range("N62").Value = -71657
$ws.Range("H65").Value = 34860.848
$ws.Range("I65").Value = 4391
$ws.Range("J65").Value = 70409
$ws.Range("K65").Value = 21955
$ws.Range("L65").Value = 352045
$ws.Range("M65").Value = -18835
$ws.Range("N65").Value = -358285
$ws.Range("H86").Value = 15092.818
$ws.Range("I86").Value = 9665
$ws.Range("K86").Value = 9665
$ws.Range("M86").Value = -8542
$ws.Range("H89").Value = 15092.818
$ws.Range("I89").Value = 9665
$ws.Range("K89").Value = 48325
$ws.Range("M89").Value = -42709
$ws.Range("H121").Value = 39999.168
$ws.Range("J121").Value = 39999.168
$ws.Range("L121").Value = 39999.168
$ws.Range("N121").Value = -42619.168
$ws.Range("H125").Value = 89666
$ws.Range("J125").Value = 89666
$ws.Range("L125").Value = 89666
$ws.Range("N125").Value = -94586
$ws.Range("H132").Value = 3319.8276
$ws.Range("I132").Value = 2855.5
$ws.Range("J132").Value = 5548.6
$ws.Range("K132").Value = 8566.5
$ws.Range("L132").Value = 16645.8
$ws.Range("M132").Value = -6036.5
$ws.Range("N132").Value = -21705.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 12887.889
$ws.Range("J32").Value = 32666.334
$ws.Range("L32").Value = 97999.00199999999
$ws.Range("N32").Value = -98565.00199999999
$ws.Range("H140").Value = 2012.7646
$ws.Range("I140").Value = 2012.7646
$ws.Range("K140").Value = 6038.293799999999
$ws.Range("M140").Value = -858.2937999999995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8736.200000000001
$ws.Range("I70").Value = 6894
$ws.Range("K70").Value = 6894
$ws.Range("M70").Value = -6624
$ws.Range("H73").Value = 8736.200000000001
$ws.Range("I73").Value = 6894
$ws.Range("K73").Value = 6894
$ws.Range("M73").Value = -5958
$ws.Range("H80").Value = 4067.375
$ws.Range("I80").Value = 2872.6667
$ws.Range("K80").Value = 2872.6667
$ws.Range("M80").Value = -1874.6667
$ws.Range("H83").Value = 4067.375
$ws.Range("I83").Value = 2872.6667
$ws.Range("K83").Value = 14363.3335
$ws.Range("M83").Value = -9371.333500000001
$ws.Range("H132").Value = 4992.6665
$ws.Range("I132").Value = 4578
$ws.Range("J132").Value = 6444
$ws.Range("K132").Value = 13734
$ws.Range("L132").Value = 19332
$ws.Range("M132").Value = -11204
$ws.Range("N132").Value = -24392

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1678.35
$ws.Range("I82").Value = 1797.5
$ws.Range("K82").Value = 1797.5
$ws.Range("M82").Value = -1436.5
$ws.Range("H85").Value = 1678.35
$ws.Range("I85").Value = 1797.5
$ws.Range("K85").Value = 1797.5
$ws.Range("M85").Value = -549.5
$ws.Range("H127").Value = 75899.8
$ws.Range("J127").Value = 75899.8
$ws.Range("L127").Value = 75899.8
$ws.Range("N127").Value = -85819.8
$ws.Range("H132").Value = 4717.64
$ws.Range("I132").Value = 3197.0908
$ws.Range("K132").Value = 9591.2724
$ws.Range("M132").Value = -7061.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 858166.7
$ws.Range("J26").Value = 2500000
$ws.Range("L26").Value = 2500000
$ws.Range("N26").Value = -2500586
$ws.Range("H60").Value = 105000
$ws.Range("J60").Value = 100000
$ws.Range("L60").Value = 100000
$ws.Range("N60").Value = -101644
$ws.Range("H100").Value = 1314.7059
$ws.Range("I100").Value = 1445.4546
$ws.Range("J100").Value = 1075
$ws.Range("K100").Value = 2890.9092
$ws.Range("L100").Value = 2150
$ws.Range("M100").Value = -2349.9092
$ws.Range("N100").Value = -3232
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

Write-Host "Updated cells: set=231 cleared=2"